# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 13583
    5  = 1031
    6  = 27
    8  = 148
    13 = 13588
    15 = 604
    16 = 8973
    18 = 8067
    21 = 151
    22 = 424
    24 = 8
    25 = 23
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
